$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.7593992856908
$ws.Range("D2").Value = 0.000353414627585295

$ws.Range("B3").Value = 7663.05747456298

$ws.Range("B4").Value = 202.757161471031
$ws.Range("D4").Value = 0.000000000000000000000000000000000000042925838767237
